$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.169.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.55%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "

# Row 6
$ws.Range("E6").Value = "  +0.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2880"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.79%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06544"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.861.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.087"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.28%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6751"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.171.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.07%  "

# Row 18
$ws.Range("E18").Value = "  +7.90%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007569"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.098.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.236"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.114"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.54%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.39%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.397"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09849"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.467"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.283"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.995"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04678"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6954"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "

# Row 37
$ws.Range("E37").Value = "  -0.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01868"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.74%  "

# Row 39
$ws.Range("E39").Value = "  +3.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.321"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.915"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "

# Row 43
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8365"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.18%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "

# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4127"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.88%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "941.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.093"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.948"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.51%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05657"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "

